$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46, pushing existing rows 46..172 down to 47..173
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new data record
$ws.Range("A46").Value = 8
$ws.Range("B46").Value = "Terminal La Palmera de La Serena"
$ws.Range("C46").Value = "Coquimbo"
$ws.Range("D46").Value = 44838
$ws.Range("E46").Value = 4
$ws.Range("F46").Value = 100112044
$ws.Range("G46").Value = "Perejil"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 2800
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 2500
$ws.Range("M46").Value = 2250
$ws.Range("N46").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O46").Value = "Provincia del Elquí"
$ws.Range("P46").Value = 1500
$ws.Range("Q46").Value = 1.5
$ws.Range("R46").Value = "Hortaliza"

# Keep the date style consistent with the other date cells in column D
$ws.Range("D46").NumberFormat = $ws.Range("D47").NumberFormat
